# Update the PCM sheet's hatch-placer / climb-actuator rows with the new
# variable names and bumped PCM-port labels (placerSol/detachLeft/detachRight
# -> scissorHolder/leftLauncher/rightLauncher, PCM port numbers shifted up by
# one, and the "HATCH/CLIMB n--m/n" sticker labels bumped accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PCM")

# Row 2: hatch pickup solenoid (PCM number 0 -> 1, sticker label bump)
$ws.Range("C2").Value = 1
$ws.Range("E2").Value = "HATCH PICKUP 1--1/6"

# Row 3: hatch placer solenoid (PCM number 1 -> 2, sticker label bump)
$ws.Range("C3").Value = 2
$ws.Range("E3").Value = "HATCH PLACER 2--1/6"

# Row 4: left hatch placing actuator (PCM number 1 -> 2, sticker label bump)
$ws.Range("C4").Value = 2
$ws.Range("E4").Value = "HATCH LEFT 2--2/5"

# Row 5: right hatch placing actuator (PCM number 1 -> 2, sticker label bump)
$ws.Range("C5").Value = 2
$ws.Range("E5").Value = "HATCH RIGHT 2--3/4"

# Row 6: left climb actuator (PCM number 0 -> 1, sticker label bump)
$ws.Range("C6").Value = 1
$ws.Range("E6").Value = "CLIMB LEFT 1--2/5"

# Row 7: right climb actuator (PCM number 0 -> 1, sticker label bump)
$ws.Range("C7").Value = 1
$ws.Range("E7").Value = "CLIMB RIGHT 1--2/5"

# Variable-name renames (placerSol/detachLeftSol/detachRightSol removed,
# replaced with scissorHolder/leftLauncher/rightLauncher)
$ws.Range("F3").Value = "scissorHolder"
$ws.Range("F4").Value = "leftLauncher"
$ws.Range("F5").Value = "rightLauncher"

# Row 8 & 9: compressor / pressure switch PCM number 0 -> 1
$ws.Range("C8").Value = 1
$ws.Range("C9").Value = 1

# Move the active selection to F6, matching the saved view state
$ws.Activate()
$ws.Range("F6").Select()
